$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: rename "Pd" -> "Pd1" and add new "Pd2" header in E1.
$ws.Range("D1").Value = "Pd1"
$ws.Range("E1").Value = "Pd2"

# Move the old demand values (previously in column D) into the new
# column E, and zero-fill column D for all data rows.
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 350

$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 300

$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 200

$ws.Range("F1").Select()
